# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share that string).
# 2) Narrow the "Status" column(s) that held that text:
#      - Overview sheet: columns E and F (zh-cn / de-de status columns)
#      - zh-cn sheet: column C (Status)
#      - de-de sheet: column C (Status)
#    from ~17.22 chars wide down to ~13.41 chars wide.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text wherever it appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Narrow the status columns ---
# Target stored width ~= 13.4101845877511 characters.
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
